# Updated symbol list refresh: new Price / Volume(1h) readings and the
# Hora (hour) column moving from 14 to 15 for every data row (rows 2-51).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Per-row Price (D) / Volume 1h (E) updates. $null means that column did not
# change for that row (e.g. rows whose Price/Volume are still "--"/"--%").
$updates = @(
    @{Row=2; D="275.62"; E="-1.20%"},
    @{Row=3; D="27.39"; E="0.99%"},
    @{Row=4; D="4.769"; E="-2.97%"},
    @{Row=5; D="0.06319"; E="-0.69%"},
    @{Row=6; D="6.928"; E="-0.20%"},
    @{Row=7; D="1.310"; E="38.36%"},
    @{Row=8; D="0.8773"; E="-0.90%"},
    @{Row=9; D="0.1522"; E="3.24%"},
    @{Row=10; D="0.05036"; E="-1.62%"},
    @{Row=11; D="0.07504"; E="0.30%"},
    @{Row=12; D="0.02874"; E="-8.43%"},
    @{Row=13; D="0.09037"; E="-0.27%"},
    @{Row=14; D="0.001570"; E="0.76%"},
    @{Row=15; D="0.0006336"; E="1.01%"},
    @{Row=16; D="0.005887"; E="-0.34%"},
    @{Row=17; D="3.450"; E="-1.03%"},
    @{Row=18; D="3.307"; E="-1.40%"},
    @{Row=19; D=$null; E="-1.10%"},
    @{Row=20; D=$null; E="0.64%"},
    @{Row=21; D="0.1325"; E="2.00%"},
    @{Row=22; D="3.896"; E="0.87%"},
    @{Row=23; D="0.04406"; E="1.68%"},
    @{Row=24; D=$null; E="-0.53%"},
    @{Row=25; D="0.003845"; E="5.73%"},
    @{Row=26; D="0.0001198"; E="-0.16%"},
    @{Row=27; D="0.0001934"; E="14.18%"},
    @{Row=28; D=$null; E=$null},
    @{Row=29; D=$null; E=$null},
    @{Row=30; D=$null; E=$null},
    @{Row=31; D=$null; E=$null},
    @{Row=32; D=$null; E=$null},
    @{Row=33; D=$null; E=$null},
    @{Row=34; D=$null; E=$null},
    @{Row=35; D=$null; E=$null},
    @{Row=36; D=$null; E=$null},
    @{Row=37; D=$null; E=$null},
    @{Row=38; D=$null; E=$null},
    @{Row=39; D=$null; E=$null},
    @{Row=40; D="0.04117"; E="1.36%"},
    @{Row=41; D="0.006929"; E="4.97%"},
    @{Row=42; D="0.1176"; E="0.84%"},
    @{Row=43; D="0.001926"; E="-18.01%"},
    @{Row=44; D=$null; E="-10.16%"},
    @{Row=45; D="0.00005213"; E="-0.81%"},
    @{Row=46; D=$null; E="-37.33%"},
    @{Row=47; D="0.01997"; E="-11.62%"},
    @{Row=48; D=$null; E=$null},
    @{Row=49; D=$null; E=$null},
    @{Row=50; D=$null; E=$null},
    @{Row=51; D=$null; E=$null}
)

foreach ($u in $updates) {
    $row = $u.Row
    if ($null -ne $u.D) {
        $ws.Range("D$row").Value = "'" + $u.D
    }
    if ($null -ne $u.E) {
        $ws.Range("E$row").Value = "'" + $u.E
    }
    # Hora (G) goes from 14 to 15 for every row in the table.
    $ws.Range("G$row").Value = "'15"
}
